$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'69.175.35"
$ws.Range("E2").Value = "'  +2.42%  "

$ws.Range("D3").Value = "'3.941.73"
$ws.Range("E3").Value = "'  +1.92%  "

$ws.Range("D4").Value = "'1.00"

$ws.Range("D5").Value = "'486.09"
$ws.Range("E5").Value = "'  +3.66%  "

$ws.Range("D6").Value = "'148.56"
$ws.Range("E6").Value = "'  -0.14%  "

$ws.Range("D7").Value = "'0.622"
$ws.Range("E7").Value = "'  -1.73%  "

$ws.Range("E8").Value = "'  -0.05%  "

$ws.Range("D9").Value = "'0.726"
$ws.Range("E9").Value = "'  -3.37%  "

$ws.Range("D10").Value = "'0.171"
$ws.Range("E10").Value = "'  +9.39%  "

$ws.Range("D11").Value = "'0.0000354"
$ws.Range("E11").Value = "'  +12.62%  "

$ws.Range("D12").Value = "'42.76"
$ws.Range("E12").Value = "'  -2.49%  "

$ws.Range("D13").Value = "'10.57"
$ws.Range("E13").Value = "'  +1.22%  "

$ws.Range("D14").Value = "'4.554.98"
$ws.Range("E14").Value = "'  +1.42%  "

$ws.Range("B15").Value = "'WrappedEther"
$ws.Range("C15").Value = "'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D15").Value = "'3.955.09"
$ws.Range("E15").Value = "'  +1.82%  "

$ws.Range("B16").Value = "'Uniswap"
$ws.Range("C16").Value = "'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D16").Value = "'14.62"
$ws.Range("E16").Value = "'  -1.24%  "

$ws.Range("E17").Value = "'  -0.23%  "

$ws.Range("D18").Value = "'19.80"
$ws.Range("E18").Value = "'  -1.49%  "

$ws.Range("E19").Value = "'  -3.27%  "

$ws.Range("D20").Value = "'69.202.60"
$ws.Range("E20").Value = "'  +2.32%  "

$ws.Range("D21").Value = "'438.44"
$ws.Range("E21").Value = "'  +1.37%  "

$ws.Range("D22").Value = "'14.64"
$ws.Range("E22").Value = "'  -1.19%  "

$ws.Range("E23").Value = "'  +0.62%  "

$ws.Range("D24").Value = "'87.45"
$ws.Range("E24").Value = "'  -1.33%  "

$ws.Range("D25").Value = "'11.56"
$ws.Range("E25").Value = "'  +14.01%  "

$ws.Range("D26").Value = "'3.58"
$ws.Range("E26").Value = "'  -0.57%  "

$ws.Range("D27").Value = "'10.58"
$ws.Range("E27").Value = "'  +2.52%  "

$ws.Range("D28").Value = "'38.27"
$ws.Range("E28").Value = "'  +1.30%  "

$ws.Range("E29").Value = "'  +6.89%  "

$ws.Range("D30").Value = "'714.01"
$ws.Range("E30").Value = "'  -4.92%  "

$ws.Range("D31").Value = "'13.27"
$ws.Range("E31").Value = "'  -3.52%  "

$ws.Range("E32").Value = "'  -4.95%  "

$ws.Range("E33").Value = "'  +2.72%  "

$ws.Range("D34").Value = "'0.0₃0905"
$ws.Range("E34").Value = "'  +33.22%  "

$ws.Range("D35").Value = "'41.46"
$ws.Range("E35").Value = "'  -3.98%  "

$ws.Range("D36").Value = "'58.52"
$ws.Range("E36").Value = "'  +1.35%  "

$ws.Range("E37").Value = "'  -7.32%  "

$ws.Range("E38").Value = "'  -0.49%  "

$ws.Range("E39").Value = "'  -0.20%  "

$ws.Range("D40").Value = "'0.0472"
$ws.Range("E40").Value = "'  -2.03%  "

$ws.Range("D41").Value = "'2.83"
$ws.Range("E41").Value = "'  +7.10%  "

$ws.Range("B42").Value = "'ThetaToken"
$ws.Range("C42").Value = "'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D42").Value = "'2.97"
$ws.Range("E42").Value = "'  +1.71%  "

$ws.Range("B43").Value = "'WEMIXToken"
$ws.Range("C43").Value = "'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("E43").Value = "'  +6.67%  "

$ws.Range("D44").Value = "'0.341"
$ws.Range("E44").Value = "'  -3.55%  "

$ws.Range("D45").Value = "'0.141"
$ws.Range("E45").Value = "'  -0.77%  "

$ws.Range("D46").Value = "'0.998"
$ws.Range("E46").Value = "'  -0.35%  "

$ws.Range("D47").Value = "'3.40"
$ws.Range("E47").Value = "'  -1.01%  "

$ws.Range("D48").Value = "'2.16"
$ws.Range("E48").Value = "'  +1.06%  "

$ws.Range("D49").Value = "'147.48"
$ws.Range("E49").Value = "'  +2.18%  "

$ws.Range("D50").Value = "'3.11"
$ws.Range("E50").Value = "'  -3.99%  "

$ws.Range("D51").Value = "'2.83"
$ws.Range("E51").Value = "'  -2.92%  "
